$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual odds values (rows 2-7) ---
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("Q2").Value = 3.4
$ws.Range("R2").Value = 1.33
$ws.Range("H3").Value = 3.3
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 2.15
$ws.Range("R3").Value = 1.67
$ws.Range("U3").Value = 1.8
$ws.Range("V3").Value = 1.95
$ws.Range("X3").Value = 11
$ws.Range("AC3").Value = 9
$ws.Range("AG3").Value = 9.5
$ws.Range("AM3").Value = 251
$ws.Range("AS3").Value = 151
$ws.Range("M4").Value = 1.13
$ws.Range("O4").Value = 1.5
$ws.Range("M5").Value = 1.11
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 2.38
$ws.Range("I6").Value = 4
$ws.Range("K6").Value = 2.05
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 8
$ws.Range("O6").Value = 1.4
$ws.Range("AC6").Value = 8
$ws.Range("AU6").Value = 8.5
$ws.Range("K7").Value = 1.91

# --- Insert a new row at position 11, shifting the old row 11 (Ameliano vs 2 de Mayo) down to row 12 ---
$ws.Rows(11).Insert()

# --- Populate the new row 11 with the Colombia Primera A match (America De Cali vs Santa Fe) ---
$ws.Range("A11").Value = "ARJPKb8t"
$ws.Range("B11").Value = "'11/11/2024"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "22:30"
$ws.Range("D11").Value = "COLOMBIA - PRIMERA A"
$ws.Range("E11").Value = "America De Cali"
$ws.Range("F11").Value = "Santa Fe"
$ws.Range("G11").Value = 1.83
$ws.Range("H11").Value = 3.1
$ws.Range("I11").Value = 5
$ws.Range("J11").Value = 2.6
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 5.5
$ws.Range("M11").Value = 1.1
$ws.Range("N11").Value = 7
$ws.Range("O11").Value = 1.44
$ws.Range("P11").Value = 2.63
$ws.Range("Q11").Value = 2.4
$ws.Range("R11").Value = 1.53
$ws.Range("S11").Value = 1.53
$ws.Range("T11").Value = 2.38
$ws.Range("U11").Value = 2.2
$ws.Range("V11").Value = 1.62
$ws.Range("W11").Value = 5.5
$ws.Range("X11").Value = 7.5
$ws.Range("Y11").Value = 9.5
$ws.Range("Z11").Value = 15
$ws.Range("AA11").Value = 19
$ws.Range("AB11").Value = 41
$ws.Range("AC11").Value = 6.5
$ws.Range("AD11").Value = 6.5
$ws.Range("AE11").Value = 19
$ws.Range("AF11").Value = 81
$ws.Range("AG11").Value = 10
$ws.Range("AH11").Value = 23
$ws.Range("AI11").Value = 17
$ws.Range("AJ11").Value = 51
$ws.Range("AK11").Value = 41
$ws.Range("AL11").Value = 51
$ws.Range("AM11").Value = 201
$ws.Range("AN11").Value = 3.6
$ws.Range("AO11").Value = 10
$ws.Range("AP11").Value = 26
$ws.Range("AQ11").Value = 41
$ws.Range("AR11").Value = 67
$ws.Range("AS11").Value = 251
$ws.Range("AT11").Value = 2.38
$ws.Range("AU11").Value = 9.5
$ws.Range("AV11").Value = 81
$ws.Range("AW11").Value = 6
$ws.Range("AX11").Value = 29
$ws.Range("AY11").Value = 41
$ws.Range("AZ11").Value = 101
$ws.Range("BA11").Value = 151
$ws.Range("BB11").Value = 351
$ws.Range("BC11").Value = 126
$ws.Range("BD11").Value = 126
